$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextCell "D2" "64.143.46"
Set-TextCell "E2" "  +1.37%  "

Set-TextCell "D3" "3.067.10"
Set-TextCell "E3" "  +0.43%  "

Set-TextCell "D4" "0.999"
Set-TextCell "E4" "  -0.11%  "

Set-TextCell "D5" "559.08"
Set-TextCell "E5" "  +1.79%  "

Set-TextCell "E6" "  +4.65%  "

Set-TextCell "E7" "  -0.02%  "

Set-TextCell "D8" "3.062.23"
Set-TextCell "E8" "  +0.49%  "

Set-TextCell "D9" "0.501"
Set-TextCell "E9" "  +0.47%  "

Set-TextCell "E10" "  +2.52%  "

Set-TextCell "D11" "6.17"

Set-TextCell "E12" "  +3.88%  "

Set-TextCell "E13" "  +0.77%  "

Set-TextCell "D14" "35.17"
Set-TextCell "E14" "  +1.56%  "

Set-TextCell "D15" "3.567.33"
Set-TextCell "E15" "  +0.29%  "

Set-TextCell "D16" "64.093.08"
Set-TextCell "E16" "  +1.18%  "

Set-TextCell "D17" "3.069.62"
Set-TextCell "E17" "  +0.29%  "

Set-TextCell "E18" "  +1.40%  "

Set-TextCell "E19" "  +1.14%  "

Set-TextCell "D20" "478.52"
Set-TextCell "E20" "  -0.16%  "

Set-TextCell "E21" "  +2.78%  "

Set-TextCell "D22" "0.676"
Set-TextCell "E22" "  +0.57%  "

Set-TextCell "D23" "7.56"
Set-TextCell "E23" "  +5.88%  "

Set-TextCell "D24" "13.56"
Set-TextCell "E24" "  +9.57%  "

Set-TextCell "D25" "81.68"
Set-TextCell "E25" "  +1.02%  "

Set-TextCell "D26" "0.998"
Set-TextCell "E26" "  -0.12%  "

Set-TextCell "D27" "2.80"
Set-TextCell "E27" "  +2.34%  "

Set-TextCell "D28" "8.10"
Set-TextCell "E28" "  +2.94%  "

Set-TextCell "E29" "  +4.90%  "

Set-TextCell "E30" "  +0.24%  "

Set-TextCell "D31" "26.22"
Set-TextCell "E31" "  +1.30%  "

Set-TextCell "E32" "  +0.75%  "

Set-TextCell "D33" "2.49"
Set-TextCell "E33" "  +3.58%  "

Set-TextCell "D34" "5.58"
Set-TextCell "E34" "  -0.93%  "

Set-TextCell "D35" "6.18"
Set-TextCell "E35" "  +3.72%  "

Set-TextCell "D36" "54.89"
Set-TextCell "E36" "  -1.29%  "

Set-TextCell "D37" "460.19"
Set-TextCell "E37" "  -0.20%  "

Set-TextCell "E38" "  +18.58%  "

Set-TextCell "E39" "  +2.33%  "

Set-TextCell "D40" "0.0405"
Set-TextCell "E40" "  +3.38%  "

Set-TextCell "D41" "2.970.45"
Set-TextCell "E41" "  -4.67%  "

Set-TextCell "E42" "  +0.65%  "

Set-TextCell "E43" "  -2.23%  "

Set-TextCell "D44" "27.89"
Set-TextCell "E44" "  -0.35%  "

Set-TextCell "D45" "0.263"
Set-TextCell "E45" "  +4.81%  "

Set-TextCell "D46" "2.14"

Set-TextCell "E47" "  +0.08%  "

Set-TextCell "E48" "  +2.57%  "

Set-TextCell "D49" "119.94"
Set-TextCell "E49" "  +3.72%  "

Set-TextCell "E50" "  +2.26%  "

Set-TextCell "D51" "2.07"
Set-TextCell "E51" "  +1.24%  "
